$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain-looking numbers as text (e.g. "30.501.38"),
# so force Text format on the data range before writing values, then restore the
# default "Normal" style so no stray number-format style is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.501.38'
$ws.Range('D3').Value = '2.110.37'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = '334.12'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').Value = '0.5264'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').Value = '0.4562'
$ws.Range('E8').Value = '  +4.89%  '
$ws.Range('D9').Value = '53.99'
$ws.Range('E9').Value = '  +15.19%  '
$ws.Range('D10').Value = '0.09019'
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').Value = '24.55'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '2.089.37'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '6.823'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '7.865'
$ws.Range('E15').Value = '  +1.71%  '
$ws.Range('D16').Value = '96.99'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = '19.48'
$ws.Range('E20').Value = '  +2.84%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').Value = '6.321'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').Value = '30.550.58'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').Value = '12.38'
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('D25').Value = '2.357'
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('D26').Value = '2.346.04'
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').Value = '2.596'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '163.74'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('D30').Value = '133.08'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '1.202'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '1.679'
$ws.Range('E33').Value = '  +9.77%  '
$ws.Range('D34').Value = '6.164'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '3.930'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').Value = '10.49'
$ws.Range('E36').Value = '  +10.30%  '
$ws.Range('D37').Value = '0.02584'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '5.602'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D39').Value = '0.06847'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').Value = '12.83'
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').Value = '0.6934'
$ws.Range('E42').Value = '  +2.29%  '
$ws.Range('D43').Value = '1.248'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = '2.379'
$ws.Range('E44').Value = '  +7.74%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').Value = '0.6409'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '14.01'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').Value = '3.656'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '0.00000000353'
$ws.Range('E49').Value = '  +24.61%  '
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('E51').Value = '  +2.07%  '

$dataRange.Style = "Normal"
